# Auto update Excel log
# Appends new sensor/alert rows to the ALERTS, Humidity and Temperature sheets.

$wb = $excel.ActiveWorkbook

function Set-RowValues {
    param($ws, $row, $date, $timestamp, $hour, $location, $value, $status)

    # Columns A ("2026-02-01") and E (e.g. "79.6%") hold date- or number-
    # shaped text. Force text formatting before assignment so Excel keeps
    # them as literal text instead of auto-converting to a date serial or a
    # percentage number, then clear the formatting override so no extra
    # style is left behind on the cell.
    $ws.Cells.Item($row, 1).NumberFormat = "@"
    $ws.Cells.Item($row, 1).Value = $date
    $ws.Cells.Item($row, 1).ClearFormats()

    $ws.Cells.Item($row, 2).Value = $timestamp
    $ws.Cells.Item($row, 3).Value = $hour
    $ws.Cells.Item($row, 4).Value = $location

    $ws.Cells.Item($row, 5).NumberFormat = "@"
    $ws.Cells.Item($row, 5).Value = $value
    $ws.Cells.Item($row, 5).ClearFormats()

    $ws.Cells.Item($row, 6).Value = $status
}

# ---------------------------------------------------------------------------
# ALERTS sheet: append rows 9-10
# ---------------------------------------------------------------------------
$alerts = $wb.Worksheets.Item("ALERTS")

Set-RowValues $alerts 9 "2026-02-01" "18:25:24" "18:00" "Bathroom" "MODERATE" "MODERATE ALERT: Bathroom occupied, no motion > 40s."
Set-RowValues $alerts 10 "2026-02-01" "18:25:40" "18:00" "Bathroom" "CRITICAL" "CRITICAL ALERT: Bathroom occupied, no motion > 60s."

# ---------------------------------------------------------------------------
# Humidity sheet: append rows 74-79
# ---------------------------------------------------------------------------
$humidity = $wb.Worksheets.Item("Humidity")

Set-RowValues $humidity 74 "2026-02-01" "18:25:24" "18:00" "Bathroom" "79.6%" "Active"
Set-RowValues $humidity 75 "2026-02-01" "18:25:28" "18:00" "Bathroom" "80.2%" "Active"
Set-RowValues $humidity 76 "2026-02-01" "18:25:33" "18:00" "Bathroom" "78.8%" "Active"
Set-RowValues $humidity 77 "2026-02-01" "18:25:38" "18:00" "Bathroom" "79.6%" "Active"
Set-RowValues $humidity 78 "2026-02-01" "18:25:43" "18:00" "Bathroom" "78.6%" "Active"
Set-RowValues $humidity 79 "2026-02-01" "18:25:48" "18:00" "Bathroom" "79.4%" "Active"

# ---------------------------------------------------------------------------
# Temperature sheet: append rows 74-79
# ---------------------------------------------------------------------------
$temperature = $wb.Worksheets.Item("Temperature")

Set-RowValues $temperature 74 "2026-02-01" "18:25:25" "18:00" "Bathroom" "29.7C" "Active"
Set-RowValues $temperature 75 "2026-02-01" "18:25:28" "18:00" "Bathroom" "29.7C" "Active"
Set-RowValues $temperature 76 "2026-02-01" "18:25:33" "18:00" "Bathroom" "29.7C" "Active"
Set-RowValues $temperature 77 "2026-02-01" "18:25:38" "18:00" "Bathroom" "29.6C" "Active"
Set-RowValues $temperature 78 "2026-02-01" "18:25:43" "18:00" "Bathroom" "29.6C" "Active"
Set-RowValues $temperature 79 "2026-02-01" "18:25:49" "18:00" "Bathroom" "29.6C" "Active"
